# Append a new "Dia 12/09: 1hr (1 dia)" line after the existing
# "Dia 09/09: 30min (1 dia)" paragraph at the end of the document body.

$d = $word.ActiveDocument

# Last paragraph currently holds "Dia 09/09: 30min (1 dia)" and sits right
# before the body's sectPr. Insert a brand-new paragraph after it; Word
# carries over the paragraph/run formatting (Arial 24, line spacing
# 360/auto, justified) from the paragraph the new one is split off from.
$lastParagraph = $d.Paragraphs.Last
$lastParagraph.Range.InsertParagraphAfter()

# The freshly inserted paragraph is now the document's last paragraph;
# give it its text.
$newParagraph = $d.Paragraphs.Last
$newParagraph.Range.Text = "Dia 12/09: 1hr (1 dia)"
